# Slide 11 (1-based) has a code-listing paragraph that reads:
#   else if(command.equals("/boardWrite.board"){
# built out of these runs:
#   "else if("  "command.equals"(err)  ("/"  "boardWrite.board"(err)  """  "){"
# The target edit turns it into:
#   else if(command.equals("/boardRegist.board"){
# where the two runs after "command.equals" merge into one run
# "command.equals("/" (no spellcheck-error flag), and the URL text becomes
# the italic run "boardRegist.board" (also without the spellcheck-error flag).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(11)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# --- Part 1: boardWrite.board -> boardRegist.board (italic), done first so it
#     does not disturb the "command.equals(" / " region handled afterwards. ---

# Insert the replacement text right after the `("/` run: the new text then
# inherits that run's (error-free) character formatting instead of the
# error-flagged formatting still sitting on "boardWrite.board".
$rSlash = $tr.Find("(`"/")
$insertStart = $rSlash.Start + $rSlash.Length
$rSlash.InsertAfter("boardRegist.board") | Out-Null

# Grab exactly the freshly inserted span and italicize it.
$newSpan = $tr.Characters($insertStart, 17)
$newSpan.Font.Italic = $true

# Remove the stale "boardWrite.board" run (now shifted right after the
# inserted text).
$rBoard = $tr.Find("boardWrite.board")
$rBoard.Text = ""

# --- Part 2: merge "command.equals" + `("/` into a single run with no
#     spellcheck-error flag. This must happen last, since re-touching a
#     spot inside an already-merged run would re-split it. ---

$rCmd = $tr.Find("command.equals")
$rCmd.Text = ""

$rSlash2 = $tr.Find("(`"/")
$rSlash2.InsertBefore("command.equals") | Out-Null
